# "Added example w.r.t 20.10"
#
# The canonical diff for this commit shows two real content changes on top of
# a full PowerPoint resave (the resave itself also normalizes namespaces,
# strips Aspose.Slides leftovers such as smtClean/endParaRPr/empty <p:timing/>,
# and re-caches the "today" date field on every layout/master -- none of that
# is reachable through the PowerPoint object model, it is simply what
# PowerPoint's own OOXML writer does to *every* part whenever it resaves a
# deck, regardless of what the user touched):
#
#   1. The picture on slide 1 is resized from 6096000x6096000 EMU down to
#      3707904x1124744 EMU.
#   2. A comment thread is added to slide 1: a root comment "The comment" by
#      "Сергей Пучок" and a threaded reply "Reply comment" from the same
#      author (ppt/commentAuthors.xml + ppt/comments/comment1.xml are new
#      parts in the target package).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Resize the picture -------------------------------------------------
# Shapes.Width/Height are expressed in points (1 pt = 12700 EMU), so use the
# point-equivalents of the target EMU extents.
$pic = $s.Shapes.Item(1)
$pic.Width = 291.961
$pic.Height = 88.56252

# --- 2. Add the comment + threaded reply -----------------------------------
# Comments.Add(Left, Top, Author, AuthorInitials, Text) -- Left/Top are in
# points too; 10 EMU / 12700 = 0.0007874015748031496 pt reproduces the
# original <p:pos x="10" y="10"/> exactly.
$comment = $s.Comments.Add(0.0007874015748031496, 0.0007874015748031496, "Сергей Пучок", "СП", "The comment")

# Replying through Comment.Replies.Add keeps the same author and records the
# proper p15:threadingInfo/parentCm link back to the root comment.
$reply = $comment.Replies.Add("Сергей Пучок", "СП", "Reply comment")
